$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need an explicit
# Text number format first, otherwise Excel auto-converts the assigned
# string into a numeric value (matching real Excel COM behavior).

$ws.Range("D2").Value = '41.456.12'
$ws.Range("E2").Value = '  -2.71%  '

$ws.Range("D3").Value = '2.464.17'
$ws.Range("E3").Value = '  -2.75%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.82'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.55'
$ws.Range("E6").Value = '  -5.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.552'
$ws.Range("E7").Value = '  -3.11%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.508'
$ws.Range("E9").Value = '  -3.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.68'
$ws.Range("E10").Value = '  -6.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0782'
$ws.Range("E11").Value = '  -2.80%  '

$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.97'
$ws.Range("E13").Value = '  -5.11%  '

$ws.Range("D14").Value = '2.842.73'
$ws.Range("E14").Value = '  -2.93%  '

$ws.Range("D15").Value = '2.460.61'
$ws.Range("E15").Value = '  -6.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.56'
$ws.Range("E16").Value = '  -8.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("E17").Value = '  -3.81%  '

$ws.Range("D18").Value = '41.461.68'
$ws.Range("E18").Value = '  -2.67%  '

$ws.Range("E19").Value = '  -6.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.58'
$ws.Range("E21").Value = '  -5.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.22'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.84'
$ws.Range("E23").Value = '  -2.10%  '

$ws.Range("E24").Value = '  -4.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.96'
$ws.Range("E25").Value = '  -4.11%  '

$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.80'
$ws.Range("E27").Value = '  -4.41%  '

$ws.Range("E28").Value = '  -3.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  -4.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.38'
$ws.Range("E30").Value = '  -7.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.26'
$ws.Range("E31").Value = '  -2.04%  '

$ws.Range("E32").Value = '  -1.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.54'
$ws.Range("E34").Value = '  -8.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0753'
$ws.Range("E35").Value = '  -5.15%  '

$ws.Range("E36").Value = '  -4.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.25'
$ws.Range("E37").Value = '  -6.30%  '

$ws.Range("E38").Value = '  -7.20%  '

$ws.Range("E39").Value = '  -5.50%  '

$ws.Range("E40").Value = '  -3.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").Value = '  -6.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.47'
$ws.Range("E42").Value = '  -1.42%  '

$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("D44").Value = '1.984.91'
$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0285'
$ws.Range("E45").Value = '  -4.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.07'
$ws.Range("E46").Value = '  -7.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.73'
$ws.Range("E47").Value = '  -1.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '77.01'
$ws.Range("E48").Value = '  -5.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.53'
$ws.Range("E49").Value = '  -4.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.58'
$ws.Range("E50").Value = '  -3.67%  '

$ws.Range("E51").Value = '  -6.28%  '
